# Apply data updates to the "Inscricoes" worksheet (Resumo de Inscricoes)
# per the commit "Data update using git".
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Inscricoes")

$ws.Range("E9").Value = 33
$ws.Range("E12").Value = 11
$ws.Range("E14").Value = 2
$ws.Range("F14").Value = 1
$ws.Range("H14").Value = 2
$ws.Range("E15").Value = 185
$ws.Range("F15").Value = 108
$ws.Range("H15").Value = 149
$ws.Range("F17").Value = 77
$ws.Range("H17").Value = 109
$ws.Range("E18").Value = 140
$ws.Range("F18").Value = 70
$ws.Range("H18").Value = 107
$ws.Range("E19").Value = 73
$ws.Range("F19").Value = 46
$ws.Range("H19").Value = 59
$ws.Range("E22").Value = 9
$ws.Range("F22").Value = 6
$ws.Range("H22").Value = 6
$ws.Range("E25").Value = 28
$ws.Range("E26").Value = 39
$ws.Range("F26").Value = 23
$ws.Range("H26").Value = 33
$ws.Range("E27").Value = 17
$ws.Range("E28").Value = 24
$ws.Range("E33").Value = 48
$ws.Range("E34").Value = 27
$ws.Range("E36").Value = 128
$ws.Range("F36").Value = 66
$ws.Range("H36").Value = 98
$ws.Range("E37").Value = 65
$ws.Range("F37").Value = 41
$ws.Range("H37").Value = 53
$ws.Range("E38").Value = 92
$ws.Range("F38").Value = 25
$ws.Range("H38").Value = 45
$ws.Range("E42").Value = 45
$ws.Range("F42").Value = 30
$ws.Range("H42").Value = 39
$ws.Range("E46").Value = 34
$ws.Range("F47").Value = 45
$ws.Range("H47").Value = 55
$ws.Range("E48").Value = 47
$ws.Range("F48").Value = 32
$ws.Range("H48").Value = 38
$ws.Range("E50").Value = 33
$ws.Range("E57").Value = 21
$ws.Range("F57").Value = 6
$ws.Range("H57").Value = 10
$ws.Range("E62").Value = 59
$ws.Range("E63").Value = 50
$ws.Range("F63").Value = 21
$ws.Range("H63").Value = 29
$ws.Range("E66").Value = 41
$ws.Range("F66").Value = 30
$ws.Range("H66").Value = 38
$ws.Range("E67").Value = 46
$ws.Range("E76").Value = 61
$ws.Range("E77").Value = 67
$ws.Range("F78").Value = 23
$ws.Range("H78").Value = 44
$ws.Range("E81").Value = 23
$ws.Range("E88").Value = 37
$ws.Range("F88").Value = 24
$ws.Range("H88").Value = 32
